# Pull new files from local to repo
# - Rename the "학생자리" (student seat) column header to "수강정원" (enrollment capacity)
# - Update the enrollment capacity value for the 인공지능 (AI) course row from 35 to 40
# - Leave the selection on the cell that was last edited (G3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G header: 학생자리 -> 수강정원
$ws.Range("G1").Value = "수강정원"

# Row 3 (인공지능 course) capacity: 35 -> 40
$ws.Range("G3").Value = 40

# Reflect the active cell selection left after the edit
$ws.Range("G3").Select()
